$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: find the first paragraph whose range text matches $needle exactly
# (ignoring the trailing paragraph mark) and whose style local name equals
# $styleName (pass $null to skip the style check), then replace that whole
# paragraph (content + end-of-paragraph mark) with the supplied OOXML.
# ---------------------------------------------------------------------------
function Set-ParagraphXml($needle, $styleName, $xml) {
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $d.Paragraphs.Item($i)
        $t = $p.Range.Text
        if ($t.Length -gt 0) {
            $t = $t.Substring(0, $t.Length - 1)
        }
        if ($t -eq $needle) {
            if ($styleName -eq $null -or $p.Style.NameLocal -eq $styleName) {
                $p.Range.InsertXML($xml)
                return $true
            }
        }
    }
    return $false
}

# ---------------------------------------------------------------------------
# 1) Title page: "Autoři:" + " " + "TýmM" (with spell-check proof markers)
#    becomes "Auto" | "r" | ":" | " " | "TýmM" as five separate runs (same
#    formatting) with the proofErr markers removed, i.e. the visible text
#    changes from "Autoři: TýmM" to "Autor: TýmM".
# ---------------------------------------------------------------------------
$authorsXml = @'
<w:p w14:paraId="61979FC7" w14:textId="07F97216" w:rsidR="00A12416" w:rsidRDefault="000529B5" w:rsidP="000C7974"><w:pPr><w:rPr><w:rStyle w:val="Zdraznnintenzivn"/><w:i w:val="0"/><w:iCs w:val="0"/><w:sz w:val="44"/><w:szCs w:val="44"/></w:rPr></w:pPr><w:r w:rsidRPr="002C004D"><w:rPr><w:rStyle w:val="Zdraznnintenzivn"/><w:i w:val="0"/><w:iCs w:val="0"/><w:sz w:val="44"/><w:szCs w:val="44"/></w:rPr><w:t>Auto</w:t></w:r><w:r w:rsidRPr="002C004D"><w:rPr><w:rStyle w:val="Zdraznnintenzivn"/><w:i w:val="0"/><w:iCs w:val="0"/><w:sz w:val="44"/><w:szCs w:val="44"/></w:rPr><w:t>r</w:t></w:r><w:r w:rsidRPr="002C004D"><w:rPr><w:rStyle w:val="Zdraznnintenzivn"/><w:i w:val="0"/><w:iCs w:val="0"/><w:sz w:val="44"/><w:szCs w:val="44"/></w:rPr><w:t>:</w:t></w:r><w:r w:rsidR="00250E28" w:rsidRPr="002C004D"><w:rPr><w:rStyle w:val="Zdraznnintenzivn"/><w:i w:val="0"/><w:iCs w:val="0"/><w:sz w:val="44"/><w:szCs w:val="44"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00B469E8" w:rsidRPr="002C004D"><w:rPr><w:rStyle w:val="Zdraznnintenzivn"/><w:i w:val="0"/><w:iCs w:val="0"/><w:sz w:val="44"/><w:szCs w:val="44"/></w:rPr><w:t>TýmM</w:t></w:r></w:p>
'@
Set-ParagraphXml "Autoři: TýmM" $null $authorsXml | Out-Null

# ---------------------------------------------------------------------------
# 2) "Přihlášení do systému ..." paragraph: the three runs that spell out
#    "Přihlášení do systému probíhá přes tlačítko " + "Prihlaseni" + ", jak
#    je zobrazeno ..." (with spell-check proof markers around "Prihlaseni")
#    are merged into a single run with no proof markers. The trailing
#    noProof " " run is untouched.
# ---------------------------------------------------------------------------
$loginText = 'Přihlášení do systému probíhá přes tlačítko Prihlaseni, jak je zobrazeno na obrázku výše. Poté je nutné zadat email a heslo a kliknout na tlačítko „Přihlásit“. Při chybném zadání hesla se vypíše příslušná hláška o špatně zadaném heslu, v případě, že email uživatele není v databázi, hláška informuje o neexistenci uživatele.'
$loginXml = '<w:p w14:paraId="07EF87F0" w14:textId="77777777" w:rsidR="001D21D3" w:rsidRDefault="001D21D3" w:rsidP="001D21D3"><w:r><w:t>' + $loginText + '</w:t></w:r><w:r><w:rPr><w:noProof/></w:rPr><w:t xml:space="preserve"> </w:t></w:r></w:p>'
Set-ParagraphXml $loginText $null $loginXml | Out-Null

# ---------------------------------------------------------------------------
# 3) "Seznam obrázků" heading: bookmark id 23 -> 22
# ---------------------------------------------------------------------------
$figuresXml = '<w:p w14:paraId="12DD29A3" w14:textId="2858748D" w:rsidR="00BD2B9B" w:rsidRPr="000470E4" w:rsidRDefault="00BD2B9B" w:rsidP="000470E4"><w:pPr><w:pStyle w:val="Nadpis1"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:ind w:left="426" w:hanging="426"/><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:bookmarkStart w:id="22" w:name="_Toc60163451"/><w:r w:rsidRPr="000470E4"><w:rPr><w:b/><w:bCs/></w:rPr><w:lastRenderedPageBreak/><w:t>Seznam obrázků</w:t></w:r><w:bookmarkEnd w:id="22"/></w:p>'
Set-ParagraphXml "Seznam obrázků" "Heading 1" $figuresXml | Out-Null

# ---------------------------------------------------------------------------
# 4) "Seznam tabulek" heading: bookmark id 24 -> 23
# ---------------------------------------------------------------------------
$tablesXml = '<w:p w14:paraId="699653B3" w14:textId="603D50AE" w:rsidR="00900F43" w:rsidRPr="000470E4" w:rsidRDefault="00900F43" w:rsidP="000470E4"><w:pPr><w:pStyle w:val="Nadpis1"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:ind w:left="426" w:hanging="426"/><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:bookmarkStart w:id="23" w:name="_Toc60163452"/><w:r w:rsidRPr="000470E4"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Seznam tabulek</w:t></w:r><w:bookmarkEnd w:id="23"/></w:p>'
Set-ParagraphXml "Seznam tabulek" "Heading 1" $tablesXml | Out-Null
